$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.419.01'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +6.29%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.719.93'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +3.54%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '331.77'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.004'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.43%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3731'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +2.74%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '48.08'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3350'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.179'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +4.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07360'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +4.08%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.33%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.351'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +5.02%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.08'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.80%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.027'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +6.41%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.722.62'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +3.71%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001066'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.85%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.61%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '82.07'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.77%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.004'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.54%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.49'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +4.76%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.090'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +2.94%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.77'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.42%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '26.401.83'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +6.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.448'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.65%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '153.12'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +3.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.369'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.71%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.382'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +16.08%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.32'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +3.83%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.915.00'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +3.72%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '130.61'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +4.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.144'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.32%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.922'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +3.09%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08601'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.77%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.693'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +3.15%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '12.65'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +3.44%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.354'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +3.67%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02318'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2150'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06172'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.42%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.449'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.51%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.223'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -4.47%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6147'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +3.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.004'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.41%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.01'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +3.82%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.895'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.23%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5939'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +5.57%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '127.33'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.90%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +4.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07181'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.68%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '76.47'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.00%  '
